$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "NA" page-number marker on row 116 moves down to the newly appended
# row 117, so row 116's "Numero de page" cell becomes blank again (matching
# the other "nothing to report" rows above it).
$ws.Range("C116").ClearContents()

# Append the new row for 2025-05-23 (same "nothing to report" entry, with the
# "NA" marker that used to sit on row 116).
$ws.Range("A117").NumberFormat = "@"
$ws.Range("A117").Value = "2025-05-23"
$ws.Range("A117").Style = $ws.Range("A116").Style

$ws.Range("B117").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("B117").Style = $ws.Range("B116").Style

$ws.Range("C117").Value = "NA"
$ws.Range("C117").Style = $ws.Range("C115").Style

$ws.Range("D117").Value = 1
$ws.Range("D117").Style = $ws.Range("D116").Style
